# fix: fixed formatting when scrapping floating point numbers
#
# 1) Three "Razon social" entries had a comma in the wrong place (typo
#    introduced by the scraper) - the comma used as a name separator
#    should have been a period.
# 2) The "Importe" column (H) was scraped using Spanish/Argentine
#    number formatting (". " thousands separator, "," decimal
#    separator) e.g. "21.560,00". Re-save them using plain
#    floating point formatting, e.g. "21560.00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the stray commas in Razon social values ---------------------
$ws.Range("E84").Value  = "FERNANDEZ. MARIO HUGO"
$ws.Range("E87").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E175").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

# --- 2) Re-format the Importe column (H2:H224) ---------------------------
# These cells hold numbers-as-text (scraped values), so we force the
# range to Text first, write the reformatted values, then restore the
# original (default) cell style so the look & feel of the sheet is
# unaffected.
$importeRange = $ws.Range("H2:H224")
$importeRange.NumberFormat = "@"

$importeFixes = @(
    @{Row=2; Value="21560.00"},
    @{Row=3; Value="5300.00"},
    @{Row=4; Value="38720.00"},
    @{Row=5; Value="104000.00"},
    @{Row=6; Value="38720.00"},
    @{Row=7; Value="340.00"},
    @{Row=8; Value="41884.00"},
    @{Row=9; Value="240352.40"},
    @{Row=10; Value="211410.00"},
    @{Row=11; Value="406.57"},
    @{Row=12; Value="313878.40"},
    @{Row=13; Value="5021.50"},
    @{Row=14; Value="382.00"},
    @{Row=15; Value="3690.00"},
    @{Row=16; Value="3598.04"},
    @{Row=17; Value="675.00"},
    @{Row=18; Value="1350.00"},
    @{Row=19; Value="360882.10"},
    @{Row=20; Value="445290.99"},
    @{Row=21; Value="10552.00"},
    @{Row=22; Value="5560.00"},
    @{Row=23; Value="68356.23"},
    @{Row=24; Value="13950.00"},
    @{Row=25; Value="10014.32"},
    @{Row=26; Value="14740.00"},
    @{Row=27; Value="3780.00"},
    @{Row=28; Value="22383.26"},
    @{Row=29; Value="730.00"},
    @{Row=30; Value="70440.00"},
    @{Row=31; Value="29796.95"},
    @{Row=32; Value="5700.00"},
    @{Row=33; Value="6300.00"},
    @{Row=34; Value="4590.00"},
    @{Row=35; Value="1250.00"},
    @{Row=36; Value="709.38"},
    @{Row=37; Value="416934.44"},
    @{Row=38; Value="840.00"},
    @{Row=39; Value="27132.24"},
    @{Row=40; Value="651.23"},
    @{Row=41; Value="148.00"},
    @{Row=42; Value="30987.40"},
    @{Row=43; Value="73572.82"},
    @{Row=44; Value="869.13"},
    @{Row=45; Value="720.00"},
    @{Row=46; Value="22830.00"},
    @{Row=47; Value="2514.98"},
    @{Row=48; Value="1087.88"},
    @{Row=49; Value="160.00"},
    @{Row=50; Value="4314.00"},
    @{Row=51; Value="4127.96"},
    @{Row=52; Value="21956.00"},
    @{Row=53; Value="1170.95"},
    @{Row=54; Value="5859.00"},
    @{Row=55; Value="25312.77"},
    @{Row=56; Value="3700.00"},
    @{Row=57; Value="15908.28"},
    @{Row=58; Value="215.25"},
    @{Row=59; Value="6470.00"},
    @{Row=60; Value="2817.37"},
    @{Row=61; Value="1150.13"},
    @{Row=62; Value="9761.93"},
    @{Row=63; Value="284800.00"},
    @{Row=64; Value="32298.26"},
    @{Row=65; Value="56251.81"},
    @{Row=66; Value="413.10"},
    @{Row=67; Value="780.00"},
    @{Row=68; Value="7912.00"},
    @{Row=69; Value="884.84"},
    @{Row=70; Value="1235.00"},
    @{Row=71; Value="919.04"},
    @{Row=72; Value="1200.00"},
    @{Row=73; Value="12800.00"},
    @{Row=74; Value="7106.71"},
    @{Row=75; Value="204846.70"},
    @{Row=76; Value="1262.92"},
    @{Row=77; Value="508.75"},
    @{Row=78; Value="830.00"},
    @{Row=79; Value="1846.00"},
    @{Row=80; Value="720.00"},
    @{Row=81; Value="11260.00"},
    @{Row=82; Value="2200.00"},
    @{Row=83; Value="7891.00"},
    @{Row=84; Value="7065.00"},
    @{Row=85; Value="75.00"},
    @{Row=86; Value="46775.00"},
    @{Row=87; Value="17360.00"},
    @{Row=88; Value="396.80"},
    @{Row=89; Value="400.00"},
    @{Row=90; Value="65.00"},
    @{Row=91; Value="469.70"},
    @{Row=92; Value="9500.00"},
    @{Row=93; Value="12000.00"},
    @{Row=94; Value="228192.15"},
    @{Row=95; Value="15120.00"},
    @{Row=96; Value="415230.93"},
    @{Row=97; Value="53261.49"},
    @{Row=98; Value="4200.00"},
    @{Row=99; Value="9.61"},
    @{Row=100; Value="96.61"},
    @{Row=101; Value="498.00"},
    @{Row=102; Value="39.00"},
    @{Row=103; Value="17156.90"},
    @{Row=104; Value="500.00"},
    @{Row=105; Value="850.00"},
    @{Row=106; Value="300.00"},
    @{Row=107; Value="571.90"},
    @{Row=108; Value="207.00"},
    @{Row=109; Value="5240.00"},
    @{Row=110; Value="304.00"},
    @{Row=111; Value="256.80"},
    @{Row=112; Value="6721.00"},
    @{Row=113; Value="4675.00"},
    @{Row=114; Value="4950.00"},
    @{Row=115; Value="63287.37"},
    @{Row=116; Value="56130.00"},
    @{Row=117; Value="4890.40"},
    @{Row=118; Value="4223.80"},
    @{Row=119; Value="49.52"},
    @{Row=120; Value="8700.00"},
    @{Row=121; Value="208.68"},
    @{Row=122; Value="6335.00"},
    @{Row=123; Value="392.00"},
    @{Row=124; Value="750.00"},
    @{Row=125; Value="2100.00"},
    @{Row=126; Value="48.00"},
    @{Row=127; Value="1538.00"},
    @{Row=128; Value="309422.00"},
    @{Row=129; Value="2748.00"},
    @{Row=130; Value="3500.00"},
    @{Row=131; Value="2000.00"},
    @{Row=132; Value="216500.00"},
    @{Row=133; Value="5850.00"},
    @{Row=134; Value="671.50"},
    @{Row=135; Value="1154.00"},
    @{Row=136; Value="7108.70"},
    @{Row=137; Value="6437.82"},
    @{Row=138; Value="2491.50"},
    @{Row=139; Value="106.40"},
    @{Row=140; Value="548.00"},
    @{Row=141; Value="501.76"},
    @{Row=142; Value="1440.76"},
    @{Row=143; Value="514800.00"},
    @{Row=144; Value="209300.00"},
    @{Row=145; Value="12963.76"},
    @{Row=146; Value="16000.00"},
    @{Row=147; Value="4700.00"},
    @{Row=148; Value="20000.00"},
    @{Row=149; Value="8000.00"},
    @{Row=150; Value="21763.06"},
    @{Row=151; Value="5000.00"},
    @{Row=152; Value="5193.00"},
    @{Row=153; Value="5112.00"},
    @{Row=154; Value="6000.00"},
    @{Row=155; Value="4000.00"},
    @{Row=156; Value="3000.00"},
    @{Row=157; Value="34645.00"},
    @{Row=158; Value="15000.00"},
    @{Row=159; Value="7000.00"},
    @{Row=160; Value="6000.00"},
    @{Row=161; Value="12523.50"},
    @{Row=162; Value="12000.00"},
    @{Row=163; Value="5000.00"},
    @{Row=164; Value="3000.00"},
    @{Row=165; Value="22800.00"},
    @{Row=166; Value="2000.00"},
    @{Row=167; Value="2176.00"},
    @{Row=168; Value="4000.00"},
    @{Row=169; Value="915.52"},
    @{Row=170; Value="49273.20"},
    @{Row=171; Value="51150.00"},
    @{Row=172; Value="135000.00"},
    @{Row=173; Value="968.00"},
    @{Row=174; Value="16.20"},
    @{Row=175; Value="2945.00"},
    @{Row=176; Value="250.00"},
    @{Row=177; Value="2209.30"},
    @{Row=178; Value="16600.00"},
    @{Row=179; Value="13280.00"},
    @{Row=180; Value="3600.00"},
    @{Row=181; Value="1001.24"},
    @{Row=182; Value="83.60"},
    @{Row=183; Value="960.00"},
    @{Row=184; Value="1051.04"},
    @{Row=185; Value="790.00"},
    @{Row=186; Value="2616.00"},
    @{Row=187; Value="320.00"},
    @{Row=188; Value="5750.00"},
    @{Row=189; Value="1500.00"},
    @{Row=190; Value="19230.00"},
    @{Row=191; Value="160.00"},
    @{Row=192; Value="3900.00"},
    @{Row=193; Value="49.20"},
    @{Row=194; Value="120.00"},
    @{Row=195; Value="176.76"},
    @{Row=196; Value="17150.00"},
    @{Row=197; Value="3320.00"},
    @{Row=198; Value="300.00"},
    @{Row=199; Value="3341.46"},
    @{Row=200; Value="43200.00"},
    @{Row=201; Value="3800.00"},
    @{Row=202; Value="4207.88"},
    @{Row=203; Value="856065.25"},
    @{Row=204; Value="2560.00"},
    @{Row=205; Value="52512.00"},
    @{Row=206; Value="237050.00"},
    @{Row=207; Value="118000.00"},
    @{Row=208; Value="203500.00"},
    @{Row=209; Value="262350.00"},
    @{Row=210; Value="247544.00"},
    @{Row=211; Value="32500.00"},
    @{Row=212; Value="27000.00"},
    @{Row=213; Value="153750.00"},
    @{Row=214; Value="406309.00"},
    @{Row=215; Value="220000.00"},
    @{Row=216; Value="167700.00"},
    @{Row=217; Value="146132.00"},
    @{Row=218; Value="110000.00"},
    @{Row=219; Value="352.11"},
    @{Row=220; Value="39480.00"},
    @{Row=221; Value="21796.00"},
    @{Row=222; Value="44800.00"},
    @{Row=223; Value="58122.50"},
    @{Row=224; Value="4260.00"}
)

foreach ($fix in $importeFixes) {
    $ws.Cells.Item($fix.Row, 8).Value = $fix.Value
}

$importeRange.Style = "Normal"
